$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2 through 9
for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 3).Value = 45221
}
